# Update cryptos list (Price and Volume(1h) columns) as produced by the
# scheduled GitHub Actions refresh job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @(newPrice, newVolume). A $null entry means the
# value is unchanged and should be left as-is.
$updates = @{
    2  = @("63.462.81", "  +1.85%  ")
    3  = @("2.478.94", "  +1.90%  ")
    4  = @("1.00", "  -0.05%  ")
    5  = @("575.70", "  +2.01%  ")
    6  = @("147.54", "  +2.00%  ")
    7  = @("0.999", "  -0.06%  ")
    8  = @("0.542", "  +1.75%  ")
    9  = @("2.482.23", "  +2.06%  ")
    10 = @($null, "  +1.95%  ")
    11 = @($null, "  +1.02%  ")
    12 = @("5.27", "  -0.10%  ")
    13 = @("0.359", "  +2.18%  ")
    14 = @("27.22", "  +2.41%  ")
    15 = @("0.0000181", "  +3.16%  ")
    16 = @("2.931.77", "  +3.10%  ")
    17 = @("63.257.10", "  +1.79%  ")
    18 = @("2.473.59", "  +1.13%  ")
    19 = @("11.57", "  +2.61%  ")
    20 = @("7.37", "  +7.38%  ")
    21 = @("329.91", "  +1.83%  ")
    22 = @("4.19", "  +1.13%  ")
    23 = @("1.97", "  +13.44%  ")
    24 = @("1.00", "  -0.04%  ")
    25 = @("66.33", "  -1.64%  ")
    26 = @("631.84", "  +14.08%  ")
    27 = @("0.0000104", "  +9.26%  ")
    28 = @("8.52", "  -1.96%  ")
    29 = @("2.589.67", $null)
    30 = @("1.51", "  +6.09%  ")
    31 = @("0.997", "  -0.18%  ")
    32 = @("8.32", "  +0.33%  ")
    33 = @("0.144", "  -2.72%  ")
    34 = @($null, "  +0.79%  ")
    35 = @("5.18", "  +7.40%  ")
    36 = @("1.51", "  -1.16%  ")
    37 = @("1.00", "  +0.10%  ")
    38 = @("0.385", "  +0.78%  ")
    39 = @("5.46", "  -2.21%  ")
    40 = @("18.84", "  +0.69%  ")
    41 = @("147.42", "  -1.86%  ")
    42 = @("1.80", "  +0.11%  ")
    43 = @("2.64", "  +13.60%  ")
    44 = @($null, "  +0.08%  ")
    45 = @("148.75", "  +0.57%  ")
    46 = @("3.74", "  +1.94%  ")
    47 = @("20.99", "  +3.41%  ")
    48 = @("0.0543", "  +1.53%  ")
    49 = @("0.609", "  +1.99%  ")
    50 = @("0.0235", "  +2.16%  ")
    51 = @("0.0924", "  -0.11%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $price = $pair[0]
    $volume = $pair[1]

    if ($null -ne $price) {
        $cell = $ws.Cells.Item($row, 4)
        # Force text storage so values like "1.00" or "0.0000181" are not
        # silently reinterpreted as numbers (matching the original file,
        # which stores these as plain strings).
        $cell.NumberFormat = "@"
        $cell.Value = $price
    }
    if ($null -ne $volume) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $volume
    }
}
